# Weekly price-list update: a new Brócoli price row for
# "Terminal Hortofrutícola Agro Chillán" is inserted as the new row 134,
# pushing every subsequent data row (old rows 134-222) down by one
# (new rows 135-223). The sheet's used range grows from A1:R222 to A1:R223.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 134 - shifts rows 134:222 down to 135:223
$ws.Rows(134).Insert()

# Populate the newly inserted row 134 with this week's entry
$ws.Cells.Item(134, 1).Value  = 7
$ws.Cells.Item(134, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value  = "Ñuble"
$ws.Cells.Item(134, 4).Value  = 44582
$ws.Cells.Item(134, 5).Value  = 16
$ws.Cells.Item(134, 6).Value  = 100112023
$ws.Cells.Item(134, 7).Value  = "Brócoli"
$ws.Cells.Item(134, 8).Value  = "Sin especificar"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 300
$ws.Cells.Item(134, 11).Value = 700
$ws.Cells.Item(134, 12).Value = 750
$ws.Cells.Item(134, 13).Value = 725
$ws.Cells.Item(134, 14).Value = "`$/unidad"
$ws.Cells.Item(134, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(134, 16).Value = 725
$ws.Cells.Item(134, 17).Value = 1
$ws.Cells.Item(134, 18).Value = "Hortaliza"
